$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.780.03'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '2.946.66'
$ws.Range("E3").Value = '  -2.35%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''570.67'
$ws.Range("E5").Value = '  -2.46%  '

$ws.Range("D6").Value = '''163.60'
$ws.Range("E6").Value = '  +1.20%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '''0.515'
$ws.Range("E8").Value = '  -0.43%  '

$ws.Range("D9").Value = '2.941.14'
$ws.Range("E9").Value = '  -2.43%  '

$ws.Range("D10").Value = '''6.65'
$ws.Range("E10").Value = '  -1.88%  '

$ws.Range("E11").Value = '  -3.78%  '

$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").Value = '''0.0000244'
$ws.Range("E13").Value = '  -4.19%  '

$ws.Range("D14").Value = '''35.04'
$ws.Range("E14").Value = '  +1.02%  '

$ws.Range("E15").Value = '  -0.48%  '

$ws.Range("D16").Value = '65.717.69'
$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").Value = '3.433.54'
$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").Value = '''7.08'
$ws.Range("E18").Value = '  +1.63%  '

$ws.Range("D19").Value = '2.942.94'
$ws.Range("E19").Value = '  -2.36%  '

$ws.Range("D20").Value = '''15.91'
$ws.Range("E20").Value = '  +13.80%  '

$ws.Range("D21").Value = '''446.36'
$ws.Range("E21").Value = '  -2.46%  '

$ws.Range("D22").Value = '''0.696'
$ws.Range("E22").Value = '  +0.95%  '

$ws.Range("D23").Value = '''7.27'
$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").Value = '''82.07'
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("D25").Value = '''2.24'
$ws.Range("E25").Value = '  -2.13%  '

$ws.Range("D26").Value = '''12.28'
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''10.03'
$ws.Range("E27").Value = '  -5.84%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.00%  '

$ws.Range("D29").Value = '''8.21'
$ws.Range("E29").Value = '  +1.88%  '

$ws.Range("D30").Value = '''2.43'
$ws.Range("E30").Value = '  +3.99%  '

$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").Value = '''0.0000101'
$ws.Range("E32").Value = '  -5.69%  '

$ws.Range("D33").Value = '''0.117'
$ws.Range("E33").Value = '  +5.43%  '

$ws.Range("D34").Value = '''27.30'
$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").Value = '''0.972'
$ws.Range("E36").Value = '  -2.47%  '

$ws.Range("D37").Value = '''5.71'
$ws.Range("E37").Value = '  -2.34%  '

$ws.Range("D38").Value = '''47.46'
$ws.Range("E38").Value = '  +9.11%  '

$ws.Range("D39").Value = '''49.14'
$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("E40").Value = '  -9.14%  '

$ws.Range("D41").Value = '''0.302'
$ws.Range("E41").Value = '  -2.58%  '

$ws.Range("E42").Value = '  -1.52%  '

$ws.Range("D43").Value = '''2.82'
$ws.Range("E43").Value = '  -6.00%  '

$ws.Range("E44").Value = '  -0.07%  '

$ws.Range("D45").Value = '''381.61'
$ws.Range("E45").Value = '  -2.15%  '

$ws.Range("D46").Value = '''0.0350'
$ws.Range("E46").Value = '  -1.65%  '

$ws.Range("D47").Value = '2.669.44'
$ws.Range("E47").Value = '  -4.73%  '

$ws.Range("E48").Value = '  -0.85%  '

$ws.Range("D50").Value = '''23.90'
$ws.Range("E50").Value = '  +0.94%  '

$ws.Range("D51").Value = '''2.16'
$ws.Range("E51").Value = '  +0.85%  '
